$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for the changed data columns (G and H are intentionally skipped: G stays constant, H is unused)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New vm_pu values per row (row index -> values for columns B,C,D,E,F,I,J,K,L,M,N)
$data = @{
    2 = @(1.02,1.057236911003858,1.072357165746389,1.064816578236727,1.078867001718659,1.058039239564317,1.062234382032436,1.075051895221895,1.067531494027487,1.081544549425495,1.024420750799954)
    3 = @(1.02,1.058322516604582,1.073163689965299,1.065755335691515,1.079821295433667,1.058388038390082,1.06297172321609,1.075675173334507,1.068285220874017,1.082316473462619,1.02467413449566)
    4 = @(1.02,1.059025246667003,1.073685671780299,1.066363286760051,1.08043927104509,1.058612574425759,1.063448506015858,1.076077925932905,1.068772815474181,1.082815815042101,1.024837802961399)
    5 = @(1.02,1.059320739229565,1.073905137621286,1.06661899141651,1.080699183076273,1.05870669135509,1.063648867097708,1.076247110865521,1.068977771782887,1.083025702985786,1.02490654008349)
    6 = @(1.02,1.059370357535171,1.073941988324233,1.066661932491836,1.080742830177371,1.058722477690248,1.063682503998433,1.07627550998742,1.069012183167416,1.08306094201112,1.024918077290318)
    7 = @(1.02,1.059029194798212,1.073688604197688,1.066366703019448,1.08044274355108,1.058613833113553,1.063451183558898,1.076080187110676,1.068775554222782,1.082818619714321,1.024838721701978)
    8 = @(1.02,1.057603740489732,1.07262971176851,1.06513372929385,1.079189409141403,1.058157357802122,1.062483637307853,1.075262648511675,1.067786243522679,1.081805454425749,1.024506442451109)
    9 = @(1.02,1.055093978101224,1.070764663134234,1.062965021309404,1.076984610847086,1.057344117164724,1.060776211603636,1.073817854162502,1.066042074015563,1.08001904296994,1.023918729856837)
    10 = @(1.02,1.053422188563526,1.069521921751972,1.061521904681817,1.075517299163169,1.056796006744734,1.05963627235844,1.072851875952035,1.064878728042778,1.078827404506792,1.023525458545823)
    11 = @(1.02,1.052698611465688,1.068983959034735,1.060897663393725,1.074882550827606,1.056557260531193,1.059142276107446,1.072432943358511,1.064374856501741,1.078311252602705,1.023354823125933)
    12 = @(1.02,1.052429890297752,1.068784159380166,1.060665888640999,1.074646868854263,1.05646836777661,1.058958724686529,1.072277234791364,1.064187676167134,1.07811950649828,1.023291389661716)
    13 = @(1.02,1.052487529752194,1.068827016002163,1.060715600749314,1.074697419273193,1.056487445154931,1.058998099767003,1.072310639250747,1.064227827885233,1.078160637781674,1.023304998699802)
    14 = @(1.02,1.052676397921221,1.068967443046203,1.060878502854198,1.074863067405698,1.056549916947618,1.059127104896995,1.072420074443424,1.06435938452344,1.078295403303281,1.023349580753215)
    15 = @(1.02,1.052792772129598,1.069053967938874,1.06097888498562,1.074965140871702,1.056588379823165,1.0592065813723,1.072487488058925,1.064440438267181,1.07837843358088,1.023377042359099)
    16 = @(1.02,1.053470216682439,1.069557627841973,1.061563346994384,1.07555943816442,1.056811821835799,1.059669048950899,1.072879665311447,1.064912165515964,1.078861656377079,1.023536775773429)
    17 = @(1.02,1.053895244882924,1.069873601744288,1.061930135623718,1.075932388157302,1.056951603521946,1.059959036974662,1.073125491989979,1.065208031595665,1.079164725503559,1.023636879712083)
    18 = @(1.02,1.054143187422073,1.070057918640342,1.062144138727984,1.076149982135687,1.057032999682513,1.060128143958577,1.073268815047817,1.065380592224186,1.079341484689023,1.023695235217444)
    19 = @(1.02,1.054227734651059,1.070120768403589,1.06221711867459,1.076224185973896,1.057060730574192,1.060185798599497,1.073317673724768,1.065439428717081,1.079401752292233,1.023715127261739)
    20 = @(1.02,1.053849640197786,1.069839699212677,1.061890776310287,1.075892368091506,1.056936620349381,1.059927927946463,1.07309912367702,1.065176289313091,1.079132210696881,1.023626142963299)
    21 = @(1.02,1.05262077964181,1.068926090122953,1.060830529583831,1.074814285619987,1.056531526406517,1.059089117735221,1.072387851218026,1.064320644922536,1.078255718883671,1.023336453880264)
    22 = @(1.02,1.051848421389189,1.068351805032724,1.060164467514195,1.07413698391486,1.056275602556469,1.058561381912512,1.071940077202597,1.063782551706572,1.077704493277985,1.023154014823464)
    23 = @(1.02,1.052257837134378,1.068656231172571,1.060517506661802,1.074495983786677,1.056411388708751,1.058841177082477,1.072177504442555,1.064067815967012,1.077996721569691,1.023250757604342)
    24 = @(1.02,1.05387024689667,1.069855018262102,1.061908560899753,1.075910451255645,1.056943391020705,1.059941984894946,1.073111038578775,1.065190632324521,1.079146902787553,1.023630994541194)
    25 = @(1.02,1.055742567192604,1.071246716477336,1.063525212858346,1.077554156943902,1.057555409238148,1.061217915017158,1.074191860896863,1.066493085202064,1.080480998423775,1.024070926399436)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}
